$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(1, 6, 11),
    @(2, 7, 12),
    @(3, 8, 13),
    @(4, 9, 14),
    @(5, 10, 15)
)

for ($r = 0; $r -lt $data.Length; $r++) {
    for ($c = 0; $c -lt $data[$r].Length; $c++) {
        $ws.Cells.Item($r + 1, $c + 1).Value = $data[$r][$c]
    }
}

$ws.Range("F7").Select() | Out-Null
